# Generate Report for Handback
#
# The handback for e2e\4a28ded3-b3a0-46a7-a25c-d047914385c2.md has now
# completed (it was previously sitting at "Ready for handoff" / showing a
# stale-handback-file error). Update the status report:
#   - Overview sheet: roll the per-file summary row from
#     "Ready for handoff" to "Handed back: in sync with en-US"
#   - zh-cn / de-de detail sheets: flip Status to the same "handed back"
#     message, stamp the new Latest Handback DateTime, and clear the old
#     "stale handback file" Error Detail now that it is resolved.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-09-03 00:52:33"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).AutoFit()

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-09-03 00:52:40"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).AutoFit()
